# Update selections and "inconsistent" flag/value cells across several
# sheets of the workbook (runmg tests + fixed small errors).

$wb = $excel.ActiveWorkbook

# pmuVoltage: just move the saved selection.
$ws1 = $wb.Worksheets.Item("pmuVoltage")
$ws1.Range("F6").Select()

# pmuCurrent: normalize row 3 values to 1, flip a few flags to 0, move selection.
$ws2 = $wb.Worksheets.Item("pmuCurrent")
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 1
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 1
$ws2.Range("J3").Value = 1
$ws2.Range("K3").Value = 1
$ws2.Range("I4").Value = 0
$ws2.Range("F7").Value = 0
$ws2.Range("I8").Value = 0
$ws2.Range("F10").Value = 0
$ws2.Range("H23").Select()

# legacyFlow: flip flags to 0, move selection.
$ws3 = $wb.Worksheets.Item("legacyFlow")
$ws3.Range("I4").Value = 0
$ws3.Range("F5").Value = 0
$ws3.Range("I6").Value = 0
$ws3.Range("F8").Value = 0
$ws3.Range("I10").Value = 0
$ws3.Range("F13").Value = 0
$ws3.Range("I15").Value = 0
$ws3.Range("F20").Value = 0
$ws3.Range("F24").Value = 0
$ws3.Range("I4").Select()

# legacyInjection: flip flags to 0, move selection.
$ws4 = $wb.Worksheets.Item("legacyInjection")
$ws4.Range("D3").Value = 0
$ws4.Range("G4").Value = 0
$ws4.Range("D5").Value = 0
$ws4.Range("G7").Value = 0
$ws4.Range("D9").Value = 0
$ws4.Range("D3").Select()

# legacyVoltage: flip flags to 0, move selection.
$ws5 = $wb.Worksheets.Item("legacyVoltage")
$ws5.Range("D3").Value = 0
$ws5.Range("D6").Value = 0
$ws5.Range("D7").Value = 0
$ws5.Range("D10").Value = 0
$ws5.Range("I13").Select()
